$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = "Prajatantra diwas"
$ws.Range("C11").Value = "Class by bhatta sir"
